$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 8.057451718567927
$ws.Cells.Item(2, 3).Value = 4.534510739735334
$ws.Cells.Item(2, 4).Value = 4.990996010574484
$ws.Cells.Item(2, 6).Value = 25.50497943890398
$ws.Cells.Item(2, 7).Value = 31.0162231883183
$ws.Cells.Item(2, 8).Value = 14.72539155171777
$ws.Cells.Item(2, 9).Value = 21.28122733489965
$ws.Cells.Item(2, 11).Value = 8.197964696691457
$ws.Cells.Item(2, 13).Value = 19.97369168425675
$ws.Cells.Item(2, 14).Value = 18.79935615229461

$ws.Cells.Item(3, 2).Value = 7.812082630667345
$ws.Cells.Item(3, 3).Value = 4.3296563454358
$ws.Cells.Item(3, 4).Value = 4.954743815330101
$ws.Cells.Item(3, 6).Value = 25.46702457895256
$ws.Cells.Item(3, 7).Value = 30.92829985372595
$ws.Cells.Item(3, 8).Value = 14.75721707724141
$ws.Cells.Item(3, 9).Value = 21.33242587985643
$ws.Cells.Item(3, 11).Value = 8.040141642406311
$ws.Cells.Item(3, 13).Value = 19.37795430612872
$ws.Cells.Item(3, 14).Value = 18.8629646206038

$ws.Cells.Item(4, 2).Value = 7.659154786328017
$ws.Cells.Item(4, 3).Value = 4.197743839832573
$ws.Cells.Item(4, 4).Value = 4.931880793628515
$ws.Cells.Item(4, 6).Value = 25.45043442579961
$ws.Cells.Item(4, 7).Value = 30.88391432375105
$ws.Cells.Item(4, 8).Value = 14.77923614996731
$ws.Cells.Item(4, 9).Value = 21.36801913582146
$ws.Cells.Item(4, 11).Value = 7.943405309268822
$ws.Cells.Item(4, 13).Value = 19.01089779209507
$ws.Cells.Item(4, 14).Value = 18.90382010386983

$ws.Cells.Item(5, 2).Value = 7.596369522077379
$ws.Cells.Item(5, 3).Value = 4.14248334733434
$ws.Cells.Item(5, 4).Value = 4.922415073246865
$ws.Cells.Item(5, 6).Value = 25.44536580849327
$ws.Cells.Item(5, 7).Value = 30.86825032818265
$ws.Cells.Item(5, 8).Value = 14.78883145486586
$ws.Cells.Item(5, 9).Value = 21.38356656995043
$ws.Cells.Item(5, 11).Value = 7.904083714310192
$ws.Cells.Item(5, 13).Value = 18.86125094934805
$ws.Cells.Item(5, 14).Value = 18.92092292377188

$ws.Cells.Item(6, 2).Value = 7.585919140772496
$ws.Cells.Item(6, 3).Value = 4.133217673397737
$ws.Cells.Item(6, 4).Value = 4.920834396661374
$ws.Cells.Item(6, 6).Value = 25.44462643743028
$ws.Cells.Item(6, 7).Value = 30.86579592216755
$ws.Cells.Item(6, 8).Value = 14.79046231270627
$ws.Cells.Item(6, 9).Value = 21.38621112357491
$ws.Cells.Item(6, 11).Value = 7.897562074039015
$ws.Cells.Item(6, 13).Value = 18.83640565415108
$ws.Cells.Item(6, 14).Value = 18.92379028672699

$ws.Cells.Item(7, 2).Value = 7.658309784942051
$ws.Cells.Item(7, 3).Value = 4.197004618614159
$ws.Cells.Item(7, 4).Value = 4.93175373430191
$ws.Cells.Item(7, 6).Value = 25.4503592142108
$ws.Cells.Item(7, 7).Value = 30.88369325096424
$ws.Cells.Item(7, 8).Value = 14.77936303682884
$ws.Cells.Item(7, 9).Value = 21.36822459472636
$ws.Cells.Item(7, 11).Value = 7.942874526515633
$ws.Cells.Item(7, 13).Value = 19.00887953836822
$ws.Cells.Item(7, 14).Value = 18.90404891902295

$ws.Cells.Item(8, 2).Value = 7.973380387083519
$ws.Cells.Item(8, 3).Value = 4.465169587809751
$ws.Cells.Item(8, 4).Value = 4.97862342307493
$ws.Cells.Item(8, 6).Value = 25.49050169724993
$ws.Cells.Item(8, 7).Value = 30.98392145324266
$ws.Cells.Item(8, 8).Value = 14.73585021802757
$ws.Cells.Item(8, 9).Value = 21.29801633743933
$ws.Cells.Item(8, 11).Value = 8.143546435968318
$ws.Cells.Item(8, 13).Value = 19.76869224689866
$ws.Cells.Item(8, 14).Value = 18.82091565569213

$ws.Cells.Item(9, 2).Value = 8.568979666077015
$ws.Cells.Item(9, 3).Value = 4.941073002070604
$ws.Cells.Item(9, 4).Value = 5.065612434566921
$ws.Cells.Item(9, 6).Value = 25.62226846368308
$ws.Cells.Item(9, 7).Value = 31.25607608947465
$ws.Cells.Item(9, 8).Value = 14.67021739178881
$ws.Cells.Item(9, 9).Value = 21.19342906957912
$ws.Cells.Item(9, 11).Value = 8.535946455478692
$ws.Cells.Item(9, 13).Value = 21.23822933987264
$ws.Cells.Item(9, 14).Value = 18.67210802011795

$ws.Cells.Item(10, 2).Value = 8.987859388450888
$ws.Cells.Item(10, 3).Value = 5.258838404966721
$ws.Cells.Item(10, 4).Value = 5.126354217088458
$ws.Cells.Item(10, 6).Value = 25.75101442751625
$ws.Cells.Item(10, 7).Value = 31.50116290153709
$ws.Cells.Item(10, 8).Value = 14.63404728369633
$ws.Cells.Item(10, 9).Value = 21.13690082810976
$ws.Cells.Item(10, 11).Value = 8.820401655394077
$ws.Cells.Item(10, 13).Value = 22.29263336019497
$ws.Cells.Item(10, 14).Value = 18.57135606181434

$ws.Cells.Item(11, 2).Value = 9.173423579601216
$ws.Cells.Item(11, 3).Value = 5.396258587759122
$ws.Cells.Item(11, 4).Value = 5.15326448541924
$ws.Cells.Item(11, 6).Value = 25.81640094316145
$ws.Cells.Item(11, 7).Value = 31.62219279504146
$ws.Cells.Item(11, 8).Value = 14.62021720813476
$ws.Cells.Item(11, 9).Value = 21.11562223792094
$ws.Cells.Item(11, 11).Value = 8.948346937392969
$ws.Cells.Item(11, 13).Value = 22.76445015907834
$ws.Cells.Item(11, 14).Value = 18.52736459283954

$ws.Cells.Item(12, 2).Value = 9.242907279386605
$ws.Cells.Item(12, 3).Value = 5.447256472553978
$ws.Cells.Item(12, 4).Value = 5.163348069580572
$ws.Cells.Item(12, 6).Value = 25.84212923587807
$ws.Cells.Item(12, 7).Value = 31.66936831213926
$ws.Cells.Item(12, 8).Value = 14.61535795903207
$ws.Cells.Item(12, 9).Value = 21.10820448501718
$ws.Cells.Item(12, 11).Value = 8.996538489311293
$ws.Cells.Item(12, 13).Value = 22.94181357880463
$ws.Cells.Item(12, 14).Value = 18.51096949678

$ws.Cells.Item(13, 2).Value = 9.227978667347184
$ws.Cells.Item(13, 3).Value = 5.436319627449833
$ws.Cells.Item(13, 4).Value = 5.16118118418909
$ws.Cells.Item(13, 6).Value = 25.83654536204821
$ws.Cells.Item(13, 7).Value = 31.65914890703623
$ws.Cells.Item(13, 8).Value = 14.61638766950577
$ws.Cells.Item(13, 9).Value = 21.1097735365024
$ws.Cells.Item(13, 11).Value = 8.986171831263727
$ws.Cells.Item(13, 13).Value = 22.903675786738
$ws.Cells.Item(13, 14).Value = 18.51448877347729

$ws.Cells.Item(14, 2).Value = 9.179156153142241
$ws.Cells.Item(14, 3).Value = 5.400475133592762
$ws.Cells.Item(14, 4).Value = 5.154096224368025
$ws.Cells.Item(14, 6).Value = 25.81849830683668
$ws.Cells.Item(14, 7).Value = 31.6260471822346
$ws.Cells.Item(14, 8).Value = 14.61980985801163
$ws.Cells.Item(14, 9).Value = 21.11499913754166
$ws.Cells.Item(14, 11).Value = 8.952317139773298
$ws.Cells.Item(14, 13).Value = 22.77906902508348
$ws.Cells.Item(14, 14).Value = 18.52601048413468

$ws.Cells.Item(15, 2).Value = 9.149146758610764
$ws.Cells.Item(15, 3).Value = 5.378383505974535
$ws.Cells.Item(15, 4).Value = 5.149742491812632
$ws.Cells.Item(15, 6).Value = 25.80756960093248
$ws.Cells.Item(15, 7).Value = 31.60594560149032
$ws.Cells.Item(15, 8).Value = 14.62195527718604
$ws.Cells.Item(15, 9).Value = 21.11828337198367
$ws.Cells.Item(15, 11).Value = 8.931545076151247
$ws.Cells.Item(15, 13).Value = 22.70256910717176
$ws.Cells.Item(15, 14).Value = 18.53310214611849

$ws.Cells.Item(16, 2).Value = 8.97562598329576
$ws.Cells.Item(16, 3).Value = 5.249712839957604
$ws.Cells.Item(16, 4).Value = 5.124580791666202
$ws.Cells.Item(16, 6).Value = 25.74687742973388
$ws.Cells.Item(16, 7).Value = 31.4934429649801
$ws.Cells.Item(16, 8).Value = 14.63500396307956
$ws.Cells.Item(16, 9).Value = 21.13838088371253
$ws.Cells.Item(16, 11).Value = 8.812006573651749
$ws.Cells.Item(16, 13).Value = 22.26162584903483
$ws.Cells.Item(16, 14).Value = 18.57426795012989

$ws.Cells.Item(17, 2).Value = 8.867848531903245
$ws.Cells.Item(17, 3).Value = 5.168939443718928
$ws.Cells.Item(17, 4).Value = 5.108957959761423
$ws.Cells.Item(17, 6).Value = 25.71138273349128
$ws.Cells.Item(17, 7).Value = 31.42685060987517
$ws.Cells.Item(17, 8).Value = 14.64368139322854
$ws.Cells.Item(17, 9).Value = 21.15184782614224
$ws.Cells.Item(17, 11).Value = 8.738264970511224
$ws.Cells.Item(17, 13).Value = 21.98898321115106
$ws.Cells.Item(17, 14).Value = 18.5999925101608

$ws.Cells.Item(18, 2).Value = 8.805392079337523
$ws.Cells.Item(18, 3).Value = 5.121810007703155
$ws.Cells.Item(18, 4).Value = 5.099904387236404
$ws.Cells.Item(18, 6).Value = 25.6916098524583
$ws.Cells.Item(18, 7).Value = 31.38944809993126
$ws.Cells.Item(18, 8).Value = 14.64891936406161
$ws.Cells.Item(18, 9).Value = 21.16001111781371
$ws.Cells.Item(18, 11).Value = 8.695716596449341
$ws.Cells.Item(18, 13).Value = 21.83143817534273
$ws.Cells.Item(18, 14).Value = 18.61496196299292

$ws.Cells.Item(19, 2).Value = 8.784167613448265
$ws.Cells.Item(19, 3).Value = 5.105738082134507
$ws.Cells.Item(19, 4).Value = 5.096827462524598
$ws.Cells.Item(19, 6).Value = 25.68502584853007
$ws.Cells.Item(19, 7).Value = 31.37693956254404
$ws.Cells.Item(19, 8).Value = 14.65073524435643
$ws.Cells.Item(19, 9).Value = 21.16284670276674
$ws.Cells.Item(19, 11).Value = 8.681288956089183
$ws.Cells.Item(19, 13).Value = 21.77797672281001
$ws.Cells.Item(19, 14).Value = 18.62006017590442

$ws.Cells.Item(20, 2).Value = 8.879370343441364
$ws.Cells.Item(20, 3).Value = 5.177607439770956
$ws.Cells.Item(20, 4).Value = 5.110628068938041
$ws.Cells.Item(20, 6).Value = 25.71509477222052
$ws.Cells.Item(20, 7).Value = 31.43384656004549
$ws.Cells.Item(20, 8).Value = 14.64273210385958
$ws.Cells.Item(20, 9).Value = 21.15037102604756
$ws.Cells.Item(20, 11).Value = 8.746129140856878
$ws.Cells.Item(20, 13).Value = 22.01808318430742
$ws.Cells.Item(20, 14).Value = 18.59723615302087

$ws.Cells.Item(21, 2).Value = 9.193518335301853
$ws.Cells.Item(21, 3).Value = 5.411031861224061
$ws.Cells.Item(21, 4).Value = 5.156180169597034
$ws.Cells.Item(21, 6).Value = 25.82377300675316
$ws.Cells.Item(21, 7).Value = 31.63573371479487
$ws.Cells.Item(21, 8).Value = 14.61879441788286
$ws.Cells.Item(21, 9).Value = 21.11344686734321
$ws.Cells.Item(21, 11).Value = 8.962268479637983
$ws.Cells.Item(21, 13).Value = 22.81570574512537
$ws.Cells.Item(21, 14).Value = 18.52261913636193

$ws.Cells.Item(22, 2).Value = 9.394227323200671
$ws.Cells.Item(22, 3).Value = 5.557520571691136
$ws.Cells.Item(22, 4).Value = 5.185327530786988
$ws.Cells.Item(22, 6).Value = 25.9004348516913
$ws.Cells.Item(22, 7).Value = 31.77549858903967
$ws.Cells.Item(22, 8).Value = 14.60535266577768
$ws.Cells.Item(22, 9).Value = 21.09304595796882
$ws.Cells.Item(22, 11).Value = 9.102000120868915
$ws.Cells.Item(22, 13).Value = 23.32933141645154
$ws.Cells.Item(22, 14).Value = 18.4753880252663

$ws.Cells.Item(23, 2).Value = 9.287547283207999
$ws.Cells.Item(23, 3).Value = 5.479896051791014
$ws.Cells.Item(23, 4).Value = 5.169829067535233
$ws.Cells.Item(23, 6).Value = 25.85900804259516
$ws.Cells.Item(23, 7).Value = 31.70019769992332
$ws.Cells.Item(23, 8).Value = 14.61232504143865
$ws.Cells.Item(23, 9).Value = 21.10359227522629
$ws.Cells.Item(23, 11).Value = 9.027578090847722
$ws.Cells.Item(23, 13).Value = 23.05595510172044
$ws.Cells.Item(23, 14).Value = 18.50045607091808

$ws.Cells.Item(24, 2).Value = 8.874162865130515
$ws.Cells.Item(24, 3).Value = 5.173690792229184
$ws.Cells.Item(24, 4).Value = 5.109873235986588
$ws.Cells.Item(24, 6).Value = 25.713414585372
$ws.Cells.Item(24, 7).Value = 31.43068094068188
$ws.Cells.Item(24, 8).Value = 14.6431605014243
$ws.Cells.Item(24, 9).Value = 21.15103737630823
$ws.Cells.Item(24, 11).Value = 8.742574223301729
$ws.Cells.Item(24, 13).Value = 22.00492956204683
$ws.Cells.Item(24, 14).Value = 18.59848174114209

$ws.Cells.Item(25, 2).Value = 8.410797454342903
$ws.Cells.Item(25, 3).Value = 4.817841378002605
$ws.Cells.Item(25, 4).Value = 5.042624136831019
$ws.Cells.Item(25, 6).Value = 25.58097823723947
$ws.Cells.Item(25, 7).Value = 31.17444346982245
$ws.Cells.Item(25, 8).Value = 14.68585951199132
$ws.Cells.Item(25, 9).Value = 21.2181651796512
$ws.Cells.Item(25, 11).Value = 8.430245562395362
$ws.Cells.Item(25, 13).Value = 20.84425464472791
$ws.Cells.Item(25, 14).Value = 18.71085176658674
